$wb = $excel.ActiveWorkbook

# --- Sheet "Raeume": rename "Mensa" -> "Mensa 1", add "Mensa 2" / "Mensa 3" ---
$wsRaeume = $wb.Worksheets.Item("Raeume")
$wsRaeume.Range("A7").Value = "Mensa 1"

$wsRaeume.Range("A11").Value = "Mensa 2"
$wsRaeume.Range("B11").Value = 1
$wsRaeume.Range("C11").Value = 1000

$wsRaeume.Range("A12").Value = "Mensa 3"
$wsRaeume.Range("B12").Value = 1
$wsRaeume.Range("C12").Value = 1000

# --- Sheet "Personal": add 5 "Aushilfskraft" rows ---
$wsPersonal = $wb.Worksheets.Item("Personal")

$wsPersonal.Range("A9").Value = "Aushilfskraft"
$wsPersonal.Range("B9").Value = 1
$wsPersonal.Range("C9").Value = "Paedagogische Fachkraft"
$wsPersonal.Range("D9").Value = "Gruppenleitung"

$wsPersonal.Range("A10").Value = "Aushilfskraft"
$wsPersonal.Range("B10").Value = 2
$wsPersonal.Range("C10").Value = "Paedagogische Fachkraft"
$wsPersonal.Range("D10").Value = "Gruppenleitung"

$wsPersonal.Range("A11").Value = "Aushilfskraft"
$wsPersonal.Range("B11").Value = 3
$wsPersonal.Range("C11").Value = "Paedagogische Fachkraft"
$wsPersonal.Range("D11").Value = "Gruppenleitung"

$wsPersonal.Range("A12").Value = "Aushilfskraft"
$wsPersonal.Range("B12").Value = 4
$wsPersonal.Range("C12").Value = "Paedagogische Fachkraft"
$wsPersonal.Range("D12").Value = "Gruppenleitung"

$wsPersonal.Range("A13").Value = "Aushilfskraft"
$wsPersonal.Range("B13").Value = 5
$wsPersonal.Range("C13").Value = "Paedagogische Fachkraft"
$wsPersonal.Range("D13").Value = "Gruppenleitung"

# --- Selections on each sheet (as last left by the author) ---
$wsRaeume.Range("C12").Select() | Out-Null
$wsPersonal.Range("G27").Select() | Out-Null

# --- Make "Personal" the active/selected tab ---
$wsPersonal.Activate() | Out-Null
